# Apply updated cryptocurrency price/volume data per Fri Apr 12 23:37:59 UTC 2024 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.918.02"
$ws.Range("E2").Value = "  -4.44%  "

$ws.Range("D3").Value = "'3.224.87"
$ws.Range("E3").Value = "  -7.92%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'592.37"
$ws.Range("E5").Value = "  -1.94%  "

$ws.Range("D6").Value = "'152.48"
$ws.Range("E6").Value = "  -11.57%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'3.215.84"
$ws.Range("E8").Value = "  -8.03%  "

$ws.Range("D9").Value = "'0.547"
$ws.Range("E9").Value = "  -10.11%  "

$ws.Range("E10").Value = "  -10.25%  "

$ws.Range("D11").Value = "'6.65"
$ws.Range("E11").Value = "  -7.82%  "

$ws.Range("D12").Value = "'0.504"
$ws.Range("E12").Value = "  -13.94%  "

$ws.Range("D13").Value = "'39.25"
$ws.Range("E13").Value = "  -14.65%  "

$ws.Range("D14").Value = "'0.0000247"
$ws.Range("E14").Value = "  -10.37%  "

$ws.Range("D15").Value = "'3.742.61"
$ws.Range("E15").Value = "  -7.94%  "

$ws.Range("D16").Value = "'66.968.99"
$ws.Range("E16").Value = "  -4.29%  "

$ws.Range("D17").Value = "'3.223.37"
$ws.Range("E17").Value = "  -8.17%  "

$ws.Range("E18").Value = "  -4.70%  "

$ws.Range("D19").Value = "'7.24"
$ws.Range("E19").Value = "  -13.27%  "

$ws.Range("D20").Value = "'533.86"
$ws.Range("E20").Value = "  -12.92%  "

$ws.Range("E21").Value = "  -13.72%  "

$ws.Range("E22").Value = "  -12.85%  "

$ws.Range("D23").Value = "'7.96"
$ws.Range("E23").Value = "  -12.83%  "

$ws.Range("D24").Value = "'13.87"
$ws.Range("E24").Value = "  -10.53%  "

$ws.Range("D25").Value = "'86.17"
$ws.Range("E25").Value = "  -12.53%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -14.17%  "

$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  -13.27%  "

$ws.Range("D29").Value = "'8.19"
$ws.Range("E29").Value = "  -8.87%  "

$ws.Range("D30").Value = "'29.51"
$ws.Range("E30").Value = "  -12.43%  "

$ws.Range("E31").Value = "  -11.45%  "

$ws.Range("E32").Value = "  -10.75%  "

$ws.Range("D33").Value = "'543.08"
$ws.Range("E33").Value = "  -13.82%  "

$ws.Range("E34").Value = "  -18.12%  "

$ws.Range("D35").Value = "'5.77"
$ws.Range("E35").Value = "  -15.11%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").Value = "'53.16"
$ws.Range("E37").Value = "  -6.18%  "

$ws.Range("D38").Value = "'0.0875"
$ws.Range("E38").Value = "  -12.14%  "

# Rows 39 and 40 swapped (VeChain now row 39, Cosmos now row 40) with updated values
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0425"
$ws.Range("E39").Value = "  -11.30%  "

$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'9.38"
$ws.Range("E40").Value = "  -12.52%  "

$ws.Range("E41").Value = "  -12.21%  "

$ws.Range("D42").Value = "'2.945.81"
$ws.Range("E42").Value = "  -12.17%  "

$ws.Range("E43").Value = "  -23.35%  "

$ws.Range("D44").Value = "'0.267"
$ws.Range("E44").Value = "  -13.59%  "

$ws.Range("D45").Value = "'0.0₃0590"
$ws.Range("E45").Value = "  -19.27%  "

$ws.Range("D46").Value = "'2.43"
$ws.Range("E46").Value = "  -16.18%  "

$ws.Range("D47").Value = "'26.56"
$ws.Range("E47").Value = "  -16.47%  "

$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("D49").Value = "'2.13"
$ws.Range("E49").Value = "  -16.40%  "

$ws.Range("E50").Value = "  -11.66%  "

$ws.Range("D51").Value = "'122.44"
$ws.Range("E51").Value = "  -8.06%  "

